# Populate "Sheet1" (1st tab, xl/worksheets/sheet1.xml) and "Sheet3"
# (2nd tab / active tab, xl/worksheets/sheet2.xml) with the cell values
# seen in the target workbook, then restore the selections shown in each
# sheet's view so the generated sharedStrings.xml / sheetData / dimension
# / selection all line up with the target OOXML.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)   # "Sheet1"  -> xl/worksheets/sheet1.xml
$ws3 = $wb.Worksheets.Item(2)   # "Sheet3"  -> xl/worksheets/sheet2.xml (tabSelected)

# Populate Sheet3 first so its strings land at shared-string indices 0-3,
# then Sheet1, then the remaining Sheet3 cells -- this reproduces the
# exact shared string table order of the target file.
$ws3.Range("C6").Value  = "诉讼时效"
$ws3.Range("E10").Value = "飒飒下啊是"
$ws3.Range("H7").Value  = "按顺序暗杀"
$ws3.Range("M9").Value  = "暗杀啊"

$ws1.Range("E12").Value = "嘻嘻嘻"
$ws1.Range("E7").Value  = "  想啊伤心啊"
$ws1.Range("B6").Value  = "暗杀暗杀"
$ws1.Range("I7").Value  = "阿斯顿撒打算"

$ws3.Range("D10").Value = "cxxx"
$ws3.Range("C10").Value = "ccc"
$ws3.Range("J15").Value = "暗杀啊"

# Restore the selections shown in the target sheetViews. Select Sheet1's
# range first, Sheet3's last, so Sheet3 ends up as the active/tabSelected
# sheet (matching activeTab="1" / tabSelected="1" in the target).
$ws1.Range("A3:K17").Select()
$ws3.Range("G11").Select()
